$d = $word.ActiveDocument

# Locate the paragraph "Abas no Chrome ficam crashando" (the bullet list item
# that needs to be split into three separate bullets).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Abas no Chrome ficam crashando*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$fragment = '<w:p w:rsidR="00953A9C" w:rsidRDefault="00953A9C" w:rsidP="009F2DC7"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Abas no </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Chrome</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> ficam </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>crashando</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Adobe Flash Player </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>crashando</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Google Drive crash ao iniciar</w:t></w:r></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($xml)
